$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 18410.059
$ws.Range("J40").Value = 13498.071
$ws.Range("L40").Value = 13498.071
$ws.Range("N40").Value = -13848.071

# Row 113
$ws.Range("H113").Value = 26337122
$ws.Range("J113").Value = 47654896
$ws.Range("L113").Value = 47654896
$ws.Range("N113").Value = -47661404

# Row 132
$ws.Range("H132").Value = 5524.42
$ws.Range("I132").Value = 1277.8125
$ws.Range("K132").Value = 3833.4375
$ws.Range("M132").Value = -1303.4375

# Row 135
$ws.Range("H135").Value = 8378.666999999999
$ws.Range("J135").Value = 12432.889
$ws.Range("L135").Value = 111896.001
$ws.Range("N135").Value = -116966.001

# Row 137
$ws.Range("H137").Value = 4459.2856
$ws.Range("I137").Value = 2129.8
$ws.Range("J137").Value = 6577
$ws.Range("K137").Value = 6389.400000000001
$ws.Range("L137").Value = 19731
$ws.Range("M137").Value = -3839.400000000001
$ws.Range("N137").Value = -24831

# Row 138
$ws.Range("H138").Value = 6774.375
$ws.Range("J138").Value = 8535.25
$ws.Range("L138").Value = 25605.75
$ws.Range("N138").Value = -35885.75

# Row 140
$ws.Range("H140").Value = 71606.19
$ws.Range("J140").Value = 69303.336
$ws.Range("L140").Value = 69303.336
$ws.Range("N140").Value = -79663.336

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2213.1667
$ws.Range("I45").Value = 2021
$ws.Range("J45").Value = 2597.5
$ws.Range("K45").Value = 2021
$ws.Range("L45").Value = 2597.5
$ws.Range("M45").Value = -1644
$ws.Range("N45").Value = -3351.5

# Row 61
$ws.Range("H61").Value = 22744688
$ws.Range("I61").Value = 35730092
$ws.Range("J61").Value = 20230.5
$ws.Range("K61").Value = 35730092
$ws.Range("L61").Value = 20230.5
$ws.Range("M61").Value = -35729880
$ws.Range("N61").Value = -20654.5

# Row 74
$ws.Range("H74").Value = 50001796
$ws.Range("I74").Value = 125001490
$ws.Range("K74").Value = 125001490
$ws.Range("M74").Value = -125000616

# Row 77
$ws.Range("H77").Value = 50001796
$ws.Range("I77").Value = 125001490
$ws.Range("K77").Value = 625007450
$ws.Range("M77").Value = -625003082

# Row 136
$ws.Range("H136").Value = 22744688
$ws.Range("I136").Value = 35730092
$ws.Range("J136").Value = 20230.5
$ws.Range("K136").Value = 107190276
$ws.Range("L136").Value = 60691.5
$ws.Range("M136").Value = -107187726
$ws.Range("N136").Value = -65791.5

# Row 140
$ws.Range("H140").Value = 115000
$ws.Range("J140").Value = 115000
$ws.Range("L140").Value = 115000
$ws.Range("N140").Value = -125360

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 21740692
$ws.Range("I107").Value = 45456610
$ws.Range("J107").Value = 1102.1666
$ws.Range("K107").Value = 45456610
$ws.Range("L107").Value = 1102.1666
$ws.Range("M107").Value = -45454690
$ws.Range("N107").Value = -4942.1666

# Row 134
$ws.Range("H134").Value = 1740.7693
$ws.Range("I134").Value = 1601.4546
$ws.Range("K134").Value = 4804.3638
$ws.Range("M134").Value = -2269.3638

# Row 140
$ws.Range("H140").Value = 243747
$ws.Range("J140").Value = 243747
$ws.Range("L140").Value = 243747
$ws.Range("N140").Value = -254107

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3796.6606
$ws.Range("I31").Value = 3007.1365
$ws.Range("J31").Value = 6691.5835
$ws.Range("K31").Value = 3007.1365
$ws.Range("L31").Value = 6691.5835
$ws.Range("M31").Value = -2712.1365
$ws.Range("N31").Value = -7281.5835

# Row 34
$ws.Range("H34").Value = 3796.6606
$ws.Range("I34").Value = 3007.1365
$ws.Range("J34").Value = 6691.5835
$ws.Range("K34").Value = 3007.1365
$ws.Range("L34").Value = 6691.5835
$ws.Range("M34").Value = -2805.1365
$ws.Range("N34").Value = -7095.5835

# Row 109
$ws.Range("H109").Value = 58758.332
$ws.Range("J109").Value = 58758.332
$ws.Range("L109").Value = 58758.332
$ws.Range("N109").Value = -60838.332

# Row 134
$ws.Range("H134").Value = 2537.0435
$ws.Range("I134").Value = 2335.8948
$ws.Range("J134").Value = 3492.5
$ws.Range("K134").Value = 7007.6844
$ws.Range("L134").Value = 10477.5
$ws.Range("M134").Value = -4472.6844
$ws.Range("N134").Value = -15547.5

# Row 141
$ws.Range("H141").Value = 85590.64999999999
$ws.Range("I141").Value = 41999.5
$ws.Range("J141").Value = 91402.8
$ws.Range("K141").Value = 41999.5
$ws.Range("L141").Value = 91402.8
$ws.Range("M141").Value = -36819.5
$ws.Range("N141").Value = -101762.8

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 37766450
$ws.Range("I4").Value = 1182126.1
$ws.Range("K4").Value = 3546378.3
$ws.Range("M4").Value = -3546266.3

$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = $null

# Row 45
$ws.Range("H45").Value = 26689.75
$ws.Range("J45").Value = 26689.75
$ws.Range("L45").Value = 26689.75
$ws.Range("N45").Value = -27807.75

# Row 113
$ws.Range("H113").Value = 1583801.8
$ws.Range("I113").Value = 3154486.8
$ws.Range("J113").Value = 13116.667
$ws.Range("K113").Value = 3154486.8
$ws.Range("L113").Value = 13116.667
$ws.Range("M113").Value = -3152316.8
$ws.Range("N113").Value = -17456.667

# Row 123
$ws.Range("H123").Value = 49970
$ws.Range("J123").Value = 49970
$ws.Range("L123").Value = 49970
$ws.Range("N123").Value = -54870

$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 80000
$ws.Range("J6").Value = 80000
$ws.Range("L6").Value = 80000
$ws.Range("N6").Value = -80224

# Row 16
$ws.Range("H16").Value = 3293.4443
$ws.Range("I16").Value = 1948.4286
$ws.Range("J16").Value = 8001
$ws.Range("K16").Value = 1948.4286
$ws.Range("L16").Value = 8001
$ws.Range("M16").Value = -1778.4286
$ws.Range("N16").Value = -8341

# Row 55
$ws.Range("H55").Value = 815.2727
$ws.Range("I55").Value = 797
$ws.Range("K55").Value = 797
$ws.Range("M55").Value = -624

# Row 63
$ws.Range("H63").Value = 49996.332
$ws.Range("J63").Value = 49996.332
$ws.Range("L63").Value = 49996.332
$ws.Range("N63").Value = -51494.332

# Row 66
$ws.Range("H66").Value = 49996.332
$ws.Range("J66").Value = 49996.332
$ws.Range("L66").Value = 149988.996
$ws.Range("N66").Value = -157476.996

# Row 123
$ws.Range("H123").Value = 59694.5
$ws.Range("J123").Value = 59694.5
$ws.Range("L123").Value = 59694.5
$ws.Range("N123").Value = -69494.5

# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null

# Row 131
$ws.Range("H131").Value = 58310.25
$ws.Range("J131").Value = 58310.25
$ws.Range("L131").Value = 58310.25
$ws.Range("N131").Value = -68390.25

# Row 139
$ws.Range("H139").Value = 85500
$ws.Range("J139").Value = 85500
$ws.Range("L139").Value = 85500
$ws.Range("N139").Value = -95780

$ws = $wb.Worksheets.Item("WVR")
# Row 102
$ws.Range("H102").Value = 79992.5
$ws.Range("J102").Value = 79992.5
$ws.Range("L102").Value = 79992.5
$ws.Range("N102").Value = -86482.5

# Row 115
$ws.Range("H115").Value = 79974.664
$ws.Range("J115").Value = 79987
$ws.Range("L115").Value = 79987
$ws.Range("N115").Value = -83121

# Row 127
$ws.Range("H127").Value = 59800
$ws.Range("J127").Value = 59800
$ws.Range("L127").Value = 59800
$ws.Range("N127").Value = -69720

# Row 132
$ws.Range("H132").Value = 17367392
$ws.Range("I132").Value = 2527147
$ws.Range("J132").Value = 50015932
$ws.Range("K132").Value = 7581441
$ws.Range("L132").Value = 150047796
$ws.Range("M132").Value = -7578911
$ws.Range("N132").Value = -150052856

# Row 136
$ws.Range("H136").Value = 9018.759
$ws.Range("I136").Value = 4010.1538
$ws.Range("K136").Value = 12030.4614
$ws.Range("M136").Value = -9480.4614
